# chinh lai selector, them log
# - B2 (supplierFileAddress) now points to a relative output path instead of
#   the hard-coded absolute UiPath project path.
# - Workbook-wide font switched from Calibri to Arial (Normal + Hyperlink
#   cell styles, which cover every cell in this sheet).
# - Row 2's fixed 43.2pt height is no longer needed once the text fits the
#   (now narrower) font, so let Excel fall back to the sheet default again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "chinh lai selector" - point at the relative Output folder instead of the
# old absolute UiPath project path.
$ws.Range("B2").Value = "Output\Danh muc NCC mau.xlsx"

# Switch the workbook's base fonts from Calibri to Arial.
$wb.Styles.Item("Normal").Font.Name = "Arial"
$wb.Styles.Item("Hyperlink").Font.Name = "Arial"

# Changing the named style's font resets any explicit cell number formats
# that rode along on it, so restore B4's text format ("@") used for the
# phone-number-like value.
$ws.Range("B4").NumberFormat = "@"

# Let row 2 re-fit to the (now updated) default row height instead of the
# old fixed 43.2pt override.
$ws.Rows.Item(2).AutoFit()
